$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding a plain count (numeric-looking text); force Text format so
# the value stays a string "N" (matching the source data) instead of being
# auto-converted to a number by Excel.
$countCells = @{
    "K2" = "3"
    "S2" = "1"
    "AE2" = "2"
    "G3" = "1"
    "O3" = "2"
    "S3" = "2"
    "C4" = "3"
    "S4" = "4"
    "W4" = "1"
    "AE7" = "1"
    "O8" = "3"
    "G9" = "1"
    "S9" = "1"
    "AA9" = "2"
    "S10" = "2"
    "AA11" = "2"
    "G12" = "2"
    "O13" = "3"
    "O14" = "3"
    "AE16" = "1"
    "O17" = "2"
    "K18" = "2"
    "C20" = "3"
    "G20" = "1"
    "G22" = "1"
    "W24" = "5"
    "AE24" = "4"
    "G26" = "1"
    "O26" = "2"
    "G27" = "2"
    "S28" = "1"
    "AE28" = "1"
    "O29" = "1"
    "G31" = "1"
    "G32" = "3"
    "G39" = "4"
    "G45" = "3"
    "G53" = "1"
}
foreach ($ref in $countCells.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $countCells[$ref]
}

# Cells holding free-form text (code lists, the update-date banner); a plain
# value assignment already round-trips as text since they are not numeric-looking.
$textCells = @{
    "L2" = "maa://39402, *maa://30515, *maa://34787"
    "T2" = "maa://22742"
    "AF2" = "maa://25251, maa://59087"
    "H3" = "maa://21247"
    "P3" = "maa://21249, maa://26254"
    "T3" = "maa://24617, maa://45854"
    "D4" = "maa://24632, maa://22499, maa://22746"
    "T4" = "maa://32509, maa://27295, maa://22754, *maa://31008"
    "X4" = "maa://43217"
    "AF7" = "maa://45272"
    "A8" = "更新日期：2025.06.10 14:12:09"
    "P8" = "maa://32931, maa://23252, maa://37496"
    "H9" = "maa://56348"
    "T9" = "maa://26222"
    "AB9" = "maa://28711, maa://40166"
    "T10" = "maa://27395, maa://22755"
    "AB11" = "maa://29912, maa://22516"
    "H12" = "maa://21867, maa://54294"
    "P13" = "maa://22676, *maa://22583, maa://48321"
    "P14" = "maa://23250, maa://20107, maa://22772"
    "AF16" = "maa://27755"
    "P17" = "maa://23890, maa://56238"
    "L18" = "maa://22466, maa://52226"
    "D20" = "maa://21432, maa://25198, maa://36680"
    "H20" = "maa://22864"
    "H22" = "maa://25236"
    "X24" = "maa://29988, maa://23504, *maa://25141, *maa://36663, maa://52227"
    "AF24" = "maa://22523, *maa://36672, maa://29910, maa://45831"
    "H26" = "maa://24913"
    "P26" = "maa://39870, maa://56625"
    "H27" = "*maa://39601, maa://34494"
    "T28" = "maa://23263"
    "AF28" = "maa://36660"
    "P29" = "maa://54169"
    "H31" = "maa://32721"
    "H32" = "maa://21895, maa://36667, maa://22760"
    "H39" = "maa://36670, maa://25199, maa://30434, *maa://45059"
    "H45" = "maa://21229, maa://30807, maa://42459"
    "H53" = "maa://32534"
}
foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}